# Update "想去人数" (interest count) figures in the generated gh-pages output.
# Workbook has 4 sheets: 展览, 演出, 本地生活, 全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 400
$ws1.Range("F4").Value  = 1119
$ws1.Range("F8").Value  = 1054
$ws1.Range("F10").Value = 316
$ws1.Range("F11").Value = 409
$ws1.Range("F15").Value = 22
$ws1.Range("F17").Value = 390
$ws1.Range("F18").Value = 432
$ws1.Range("F19").Value = 5484
$ws1.Range("F21").Value = 1534
$ws1.Range("F22").Value = 354
$ws1.Range("F27").Value = 1468
$ws1.Range("F28").Value = 10
$ws1.Range("F31").Value = 22

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 102

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2112

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2112
$ws4.Range("F6").Value  = 400
$ws4.Range("F7").Value  = 1119
$ws4.Range("F11").Value = 1054
$ws4.Range("F13").Value = 316
$ws4.Range("F14").Value = 409
$ws4.Range("F18").Value = 22
$ws4.Range("F23").Value = 390
$ws4.Range("F24").Value = 432
$ws4.Range("F25").Value = 5484
$ws4.Range("F27").Value = 1468
$ws4.Range("F30").Value = 354
$ws4.Range("F36").Value = 1468
$ws4.Range("F37").Value = 10
$ws4.Range("F40").Value = 22
